# Insert a new weekly price-report row at row 458 (shifts existing rows
# 458:518 down to 459:519) and populate it with the new record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 458, pushing the rest down.
$ws.Rows.Item(458).Insert()

# Populate the newly inserted row with the new data record.
$ws.Cells.Item(458, 1).Value = 3
$ws.Cells.Item(458, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(458, 3).Value = "Coquimbo"
$ws.Cells.Item(458, 4).Value = 44984
$ws.Cells.Item(458, 5).Value = 5
$ws.Cells.Item(458, 6).Value = 100112040
$ws.Cells.Item(458, 7).Value = "Cilantro"
$ws.Cells.Item(458, 8).Value = "Sin especificar"
$ws.Cells.Item(458, 9).Value = "Primera"
$ws.Cells.Item(458, 10).Value = 135
$ws.Cells.Item(458, 11).Value = 5500
$ws.Cells.Item(458, 12).Value = 6000
$ws.Cells.Item(458, 13).Value = 5741
$ws.Cells.Item(458, 14).Value = "`$/docena de atados (3 kilos)"
$ws.Cells.Item(458, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(458, 16).Value = 1914
$ws.Cells.Item(458, 17).Value = 3
$ws.Cells.Item(458, 18).Value = "Hortaliza"
